$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $rng = $ws.Range($addr)
    # Force text interpretation so numeric-looking / date-looking strings
    # (ids, amounts, ISO dates) are stored verbatim instead of being
    # coerced into numbers / date serials by Excel's smart input parsing.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    # Drop the temporary formatting again so the cell keeps the sheet's
    # default (unstyled) look, matching a plain data row.
    $rng.ClearFormats()
}

Set-TextValue "A2" "ocds-twb234-0005"
Set-TextValue "B2" "3568999"
Set-TextValue "C2" "Activo"
Set-TextValue "D2" "2019-03-16"
Set-TextValue "E2" "30628707093"
Set-TextValue "F2" "HAL2000"
Set-TextValue "G2" "ARS"
Set-TextValue "H2" "2000000"
